$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 15.28031864006898
    "C2" = 10.92053664895577
    "D2" = 5.060179375518758
    "E2" = 12.55291291216698
    "F2" = 25.98620036343338
    "H2" = 7.344005520526261
    "L2" = 9.860504489982572
    "M2" = 14.82307386167018
    "N2" = 18.24257211158297
    "O2" = 23.16192675530845
    "B3" = 14.8409310878603
    "C3" = 10.75685622766334
    "D3" = 5.032157787610059
    "E3" = 12.59217483257692
    "F3" = 25.95894039593519
    "H3" = 7.344005520526261
    "L3" = 9.867969407650083
    "M3" = 14.73780870722297
    "N3" = 18.30153048265199
    "O3" = 23.19811678039135
    "B4" = 14.56694563306081
    "C4" = 10.65417043017128
    "D4" = 5.014649981082602
    "E4" = 12.61762104840495
    "F4" = 25.94989977007672
    "H4" = 7.344005520526261
    "L4" = 9.87390073256962
    "M4" = 14.68755931435047
    "N4" = 18.33958109138847
    "O4" = 23.22626614977677
    "B5" = 14.45441237404868
    "C5" = 10.61180508330514
    "D5" = 5.007441558836981
    "E5" = 12.62832814387292
    "F5" = 25.9481535783078
    "H5" = 7.344005520526261
    "L5" = 9.876657197928946
    "M5" = 14.66762743703641
    "N5" = 18.35555340575668
    "O5" = 23.23922446858817
    "B6" = 14.43567834406963
    "C6" = 10.60473981807145
    "D6" = 5.006240209850445
    "E6" = 12.6301264598742
    "F6" = 25.94798070707379
    "H6" = 7.344005520526261
    "L6" = 9.877135418629708
    "M6" = 14.66435114725329
    "N6" = 18.35823380028433
    "O6" = 23.24146589712074
    "B7" = 14.5654313031822
    "C7" = 10.65360114238072
    "D7" = 5.014553061891507
    "E7" = 12.61776408022723
    "F7" = 25.94986837183668
    "H7" = 7.344005520526261
    "L7" = 9.873936532414767
    "M7" = 14.68728827799375
    "N7" = 18.33979460964714
    "O7" = 23.22643489363581
    "B8" = 15.12979048646475
    "C8" = 10.86457061348766
    "D8" = 5.050581620160393
    "E8" = 12.56617294193402
    "F8" = 25.97520559019661
    "H8" = 7.344005520526261
    "L8" = 9.862798955286916
    "M8" = 14.79324771467097
    "N8" = 18.26251772875179
    "O8" = 23.17317262950706
    "B9" = 16.19609998333192
    "C9" = 11.25969596024838
    "D9" = 5.118738642219849
    "E9" = 12.47559282237654
    "F9" = 26.08577796048751
    "H9" = 7.344005520526261
    "L9" = 9.851629899209073
    "M9" = 15.0169840912392
    "N9" = 18.12560128150691
    "O9" = 23.11590196647014
    "B10" = 16.94625651222396
    "C10" = 11.53707234610262
    "D10" = 5.167167530514183
    "E10" = 12.41544858362903
    "F10" = 26.20376620621704
    "H10" = 7.344005520526261
    "L10" = 9.84989650544872
    "M10" = 15.1900504917274
    "N10" = 18.03384609527445
    "O10" = 23.10273681944917
    "B11" = 17.27871705337092
    "C11" = 11.66013701608433
    "D11" = 5.188815307826547
    "E11" = 12.38946742029901
    "F11" = 26.26530547211679
    "H11" = 7.344005520526261
    "L11" = 9.850504878624488
    "M11" = 15.27043723144011
    "N11" = 17.99400693772252
    "O11" = 23.10304496856797
    "B12" = 17.40323155515519
    "C12" = 11.7062655708246
    "D12" = 5.196955464099919
    "E12" = 12.37982649614574
    "F12" = 26.28972708571114
    "H12" = 7.344005520526261
    "L12" = 9.850935257959623
    "M12" = 15.30109579543413
    "N12" = 17.97919302236478
    "O12" = 23.1040675806908
    "B13" = 17.37647829724849
    "C13" = 11.69635240513993
    "D13" = 5.195204925244576
    "E13" = 12.38189406510102
    "F13" = 26.28441795668674
    "H13" = 7.344005520526261
    "L13" = 9.850833687665066
    "M13" = 15.2944835668582
    "N13" = 17.98237136901756
    "O13" = 23.10380705254712
    "B14" = 17.288989245736
    "C14" = 11.66394166794707
    "D14" = 5.189486166871999
    "E14" = 12.38867029978468
    "F14" = 26.26729232720144
    "H14" = 7.344005520526261
    "L14" = 9.850536283134289
    "M14" = 15.27295528279243
    "N14" = 17.99278273625355
    "O14" = 23.1031109435508
    "B15" = 17.23521653595523
    "C15" = 11.64402677736154
    "D15" = 5.185975719458937
    "E15" = 12.39284664959581
    "F15" = 26.25694756666071
    "H15" = 7.344005520526261
    "L15" = 9.850380133170189
    "M15" = 15.25979634208012
    "N15" = 17.99919543220457
    "O15" = 23.10280253496509
    "B16" = 16.92434228651568
    "C16" = 11.52896489405532
    "D16" = 5.165744894143518
    "E16" = 12.41717419906017
    "F16" = 26.19990167506925
    "H16" = 7.344005520526261
    "L16" = 9.849884783420828
    "M16" = 15.18482864140859
    "N16" = 18.03648783840165
    "O16" = 23.10284344524692
    "B17" = 16.73129274875418
    "C17" = 11.45756179810287
    "D17" = 5.153234227316588
    "E17" = 12.43245100431882
    "F17" = 26.1669117708592
    "H17" = 7.344005520526261
    "L17" = 9.84993808799614
    "M17" = 15.13924850643928
    "N17" = 18.05985161768624
    "O17" = 23.10448188293389
    "B18" = 16.61943471032325
    "C18" = 11.41620084453267
    "D18" = 5.146002534002219
    "E18" = 12.44136763368883
    "F18" = 26.14867811681931
    "H18" = 7.344005520526261
    "L18" = 9.850100287624159
    "M18" = 15.11318932803741
    "N18" = 18.07346881128057
    "O18" = 23.10601699523918
    "B19" = 16.58142424528977
    "C19" = 11.40214741711145
    "D19" = 5.143547914787119
    "E19" = 12.44440896805168
    "F19" = 26.14263219666323
    "H19" = 7.344005520526261
    "L19" = 9.850177823040829
    "M19" = 15.10439377719947
    "N19" = 18.07811012784884
    "O19" = 23.10663853512581
    "B20" = 16.751929062773
    "C20" = 11.46519316581733
    "D20" = 5.154569740545899
    "E20" = 12.43081133193897
    "F20" = 26.17034697311228
    "H20" = 7.344005520526261
    "L20" = 9.849918805771637
    "M20" = 15.14408445657705
    "N20" = 18.05734598862032
    "O20" = 23.10424611824244
    "B21" = 17.31472526009178
    "C21" = 11.67347453021116
    "D21" = 5.191167480896601
    "E21" = 12.38667459858298
    "F21" = 26.27229230972784
    "H21" = 7.344005520526261
    "L21" = 9.850618217406247
    "M21" = 15.2792729146954
    "N21" = 17.98971728291958
    "O21" = 23.10329082131172
    "B22" = 17.67444871893896
    "C22" = 11.80682842769731
    "D22" = 5.214750650228772
    "E22" = 12.35898001991695
    "F22" = 26.34542889778546
    "H22" = 7.344005520526261
    "L22" = 9.852240625375268
    "M22" = 15.36888644649273
    "N22" = 17.94710478287547
    "O22" = 23.1079465867072
    "B23" = 17.48323412127833
    "C23" = 11.73591646108008
    "D23" = 5.202195343185982
    "E23" = 12.37365601483714
    "F23" = 26.30580368864528
    "H23" = 7.344005520526261
    "L23" = 9.851268398918403
    "M23" = 15.32094959183533
    "N23" = 17.96970301671952
    "O23" = 23.10497863107965
    "B24" = 16.74260209002853
    "C24" = 11.46174398927726
    "D24" = 5.153966076754857
    "E24" = 12.43155221102696
    "F24" = 26.16879163492865
    "H24" = 7.344005520526261
    "L24" = 9.84992711342406
    "M24" = 15.14189766888974
    "N24" = 18.05847820701938
    "O24" = 23.10435085983844
    "B25" = 15.9128984735117
    "C25" = 11.15496236602522
    "D25" = 5.100580398536101
    "E25" = 12.49896871819296
    "F25" = 26.04937746147574
    "H25" = 7.344005520526261
    "L25" = 9.853512081559916
    "M25" = 14.95485454697824
    "N25" = 18.16108336179713
    "O25" = 23.12632647971817
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
